$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark column D (Price) as text so numeric-looking values
# like "1.003" or "277.29" are kept as plain text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '20.156.68'
$ws.Range('E2').Value = '  -1.13%  '

$ws.Range('D3').Value = '1.427.82'
$ws.Range('E3').Value = '  -0.88%  '

$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '0.9958'
$ws.Range('E5').Value = '  -0.58%  '

$ws.Range('D6').Value = '277.29'
$ws.Range('E6').Value = '  -0.24%  '

$ws.Range('D7').Value = '0.3701'
$ws.Range('E7').Value = '  -0.70%  '

$ws.Range('D8').Value = '0.3152'
$ws.Range('E8').Value = '  +1.64%  '

$ws.Range('D9').Value = '40.44'
$ws.Range('E9').Value = '  -0.64%  '

$ws.Range('D10').Value = '1.057'
$ws.Range('E10').Value = '  +3.92%  '

$ws.Range('D11').Value = '0.06584'
$ws.Range('E11').Value = '  -0.33%  '

$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.12%  '

$ws.Range('D13').Value = '5.550'
$ws.Range('E13').Value = '  +3.00%  '

$ws.Range('D14').Value = '18.20'
$ws.Range('E14').Value = '  +4.83%  '

$ws.Range('D15').Value = '6.219'
$ws.Range('E15').Value = '  +0.82%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.432.34'
$ws.Range('E16').Value = '  -0.47%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.00001029'
$ws.Range('E17').Value = '  +1.91%  '

$ws.Range('D18').Value = '0.05744'
$ws.Range('E18').Value = '  -1.99%  '

$ws.Range('D19').Value = '0.9945'
$ws.Range('E19').Value = '  -0.71%  '

$ws.Range('D20').Value = '71.67'
$ws.Range('E20').Value = '  -6.68%  '

$ws.Range('D21').Value = '5.625'
$ws.Range('E21').Value = '  -2.18%  '

$ws.Range('D22').Value = '14.83'
$ws.Range('E22').Value = '  +2.71%  '

$ws.Range('D23').Value = '11.14'
$ws.Range('E23').Value = '  +1.23%  '

$ws.Range('D24').Value = '2.237'
$ws.Range('E24').Value = '  -3.71%  '

$ws.Range('D25').Value = '20.213.83'
$ws.Range('E25').Value = '  -0.80%  '

$ws.Range('D26').Value = '2.313'
$ws.Range('E26').Value = '  +1.40%  '

$ws.Range('D27').Value = '135.46'
$ws.Range('E27').Value = '  -4.87%  '

$ws.Range('D28').Value = '17.46'
$ws.Range('E28').Value = '  +2.04%  '

$ws.Range('D29').Value = '1.593.65'
$ws.Range('E29').Value = '  -0.47%  '

$ws.Range('D30').Value = '111.83'
$ws.Range('E30').Value = '  +1.32%  '

$ws.Range('D31').Value = '3.948'
$ws.Range('E31').Value = '  -0.24%  '

$ws.Range('D32').Value = '5.316'
$ws.Range('E32').Value = '  -3.32%  '

$ws.Range('D33').Value = '0.8409'
$ws.Range('E33').Value = '  -9.69%  '

$ws.Range('D34').Value = '0.07802'
$ws.Range('E34').Value = '  +0.89%  '

$ws.Range('D35').Value = '1.496'
$ws.Range('E35').Value = '  +11.48%  '

$ws.Range('D36').Value = '0.05916'
$ws.Range('E36').Value = '  +2.80%  '

$ws.Range('D37').Value = '4.928'
$ws.Range('E37').Value = '  +3.20%  '

$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Value = '0.9951'
$ws.Range('E38').Value = '  -0.59%  '

$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '10.78'
$ws.Range('E39').Value = '  -2.99%  '

$ws.Range('D40').Value = '7.831'
$ws.Range('E40').Value = '  -6.74%  '

$ws.Range('D41').Value = '0.02066'
$ws.Range('E41').Value = '  +1.62%  '

$ws.Range('D42').Value = '1.111'
$ws.Range('E42').Value = '  -1.96%  '

$ws.Range('D43').Value = '0.1887'
$ws.Range('E43').Value = '  -2.53%  '

$ws.Range('D44').Value = '0.5369'
$ws.Range('E44').Value = '  +0.23%  '

$ws.Range('D45').Value = '12.40'
$ws.Range('E45').Value = '  +2.41%  '

$ws.Range('D46').Value = '3.560'
$ws.Range('E46').Value = '  -0.96%  '

$ws.Range('D47').Value = '120.19'
$ws.Range('E47').Value = '  +7.05%  '

$ws.Range('D48').Value = '0.5277'
$ws.Range('E48').Value = '  +1.71%  '

$ws.Range('D49').Value = '1.801'
$ws.Range('E49').Value = '  +0.57%  '

$ws.Range('D50').Value = '1.045'
$ws.Range('E50').Value = '  -1.30%  '

$ws.Range('D51').Value = '0.06278'
$ws.Range('E51').Value = '  +0.04%  '

# Restore the default cell style on column D so no visible formatting changes.
$ws.Range("D2:D51").Style = "Normal"
